$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New gyroscope data (x,y,z) covering rows 2:31 after the May 9th update
$data = @(
    @(0.4665061936658925, -0.1984146457969052, 0.8845337243641089),
    @(0.7887173727446912, -0.2283256008636717, -2.244795449808533),
    @(2.415960978059204, -0.9742195974378011, -3.041498450671924),
    @(-5.063323306102351, -1.817073576590591, 2.776077873566551),
    @(-6.820184679592355, 0.3425412797460279, 5.666934162962662),
    @(-0.2362654770122372, 5.465962034814462, 3.005506543552197),
    @(4.696316433887858, -2.277570799285295, -2.875306896134922),
    @(6.539746003992415, -0.2706960135815262, -1.030937138725735),
    @(4.692681789398195, 1.697400997666751, -1.396730149493496),
    @(-4.347883439531461, -0.5452939552419174, 0.9300998636320361),
    @(-9.689171734978153, -3.703270496106596, 9.968994430467122),
    @(2.393433683058804, -9.479072037865135, 1.899762265822456),
    @(4.808085488338101, -6.932880869098746, -1.150131797089299),
    @(6.954113165537517, 4.074747403462728, -1.481094896793366),
    @(0.2588006748873219, 3.881287055857073, 1.234817659153635),
    @(-3.794551531473794, -1.161946320066274, 2.955449520372877),
    @(-3.552903596092673, 6.386371462952866, 1.291433391325612),
    @(-2.714477875653444, 9.294961256139448, -3.878858893233145),
    @(2.263628258424627, -0.6571153846441176, -3.716604506268235),
    @(3.927821926042158, -5.819554076475262, -0.2955599602531069),
    @(2.031181216239927, -3.260402391938594, 0.723365492680494),
    @(-5.857060266476065, -0.9245907862981151, 6.962034623996895),
    @(-3.386064576167785, -0.2334743749861628, 3.267812457739236),
    @(3.026514137492571, -1.510893523693044, -3.966962477740206),
    @(3.048832594179637, -2.553463636660109, -2.418893182978909),
    @(2.535124498255125, -1.166468461500665, -2.345369084208624),
    @(-3.627002197153437, -0.152242738972691, 0.2601947556523713),
    @(-3.252175232943356, 1.719280733197323, 0.07690661035332244),
    @(-2.007904552945881, 8.391763182247416, -2.262497771020045),
    @(2.766276373582728, 6.427707840414632, -4.763592523687029)
)

$numRows = $data.Count
$numCols = 3

$values = New-Object 'object[,]' $numRows,$numCols
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt $numCols; $j++) {
        $values[$i, $j] = $data[$i][$j]
    }
}

$startRow = 2
$endRow = $startRow + $numRows - 1
$targetRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 3))
$targetRange.Value = $values
